# Add the new email address to the "Emails" sheet and make it the
# active/selected sheet (matching the Excel session state captured in the
# commit), with B23 left as the active cell selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Emails")

# New row with the additional e-mail address.
$ws.Range("A2").Value = "mu@wapi.com"

# Make "Emails" the active sheet (tabSelected / activeTab) and restore the
# cell selection that was captured for this sheet.
$ws.Activate()
$ws.Range("B23").Select()
